$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = '[16, 26, 119, 130, 230, 247, 264, 267, 272, 351, 390, 552, 567, 617, 724, 736, 770, 807, 851, 860]'
$ws.Range("I2").Value = '{''Numk'': 12, ''KPar'': 11, ''Bucket_index'': 500}'
$ws.Range("J2").Value = 10.13569731893949
$ws.Range("O2").Value = '[224, 225, 628, 634, 641, 905, 963, 964, 965, 966, 967]'
$ws.Range("P2").Value = 0.6885245901639345
$ws.Range("Q2").Value = '{''window_size'': 223.0, ''max_size_ae'': 3}'
$ws.Range("R2").Value = 0.5312427189201117
$ws.Range("S2").Value = '[500, 963, 964, 965, 966]'
$ws.Range("T2").Value = 0.3809523809523809
$ws.Range("U2").Value = '{''max_features'': 7, ''window_size'': 225, ''n_estimator'': 24}'
$ws.Range("V2").Value = 37.9295930180233

$ws.Range("G3").Value = '[0, 1, 2, 6, 7, 8, 9, 10, 11, 12, 13, 14, 16, 17, 20, 22, 23, 29, 30, 32, 33, 39, 61, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 76, 77, 80, 81, 82, 83, 84, 85, 87, 88, 89, 90, 91, 93, 97, 98, 102, 104, 105, 106, 108, 109, 111, 114, 119, 121, 122, 125, 126, 128, 132, 133, 136, 137, 138, 139, 140, 141, 144, 146, 147, 148, 149, 151, 152, 153, 154, 155, 156, 160, 163, 164, 165, 167, 170, 171, 172, 173, 175, 176, 177, 178, 179, 181, 182, 183, 184, 185, 186, 187, 188, 191, 193, 195, 197, 200, 201, 202, 203, 205, 206, 207, 209, 211, 213, 215, 217, 218, 219, 220, 221, 222, 223, 225, 226, 227, 228, 230, 231, 233, 235, 236, 237, 239, 240, 242, 247, 248, 250, 251, 252, 253, 254, 255, 256, 257, 259, 260, 261, 265, 267, 268, 270, 272, 273, 274, 279, 280, 282, 283, 284, 285, 286, 287, 288, 289, 290, 291, 292, 293, 294, 295, 296, 297, 300, 301, 302, 303, 304, 305, 306, 308, 309, 310, 311, 312, 313, 316, 317, 318, 319, 320, 321, 322, 323, 324, 328, 329, 330, 331, 332, 333, 335, 337, 338, 340, 341, 343, 344, 345, 348, 350, 352, 353, 355, 356, 361, 362, 363, 365, 366, 367, 368, 369, 371, 374, 375, 378, 379, 380, 381, 382, 384, 385, 386, 387, 388, 389, 391, 392, 394, 395, 396, 397, 398, 399, 400, 401, 402, 403, 405, 406, 408, 409, 410, 411, 412, 413, 414, 415, 416, 417, 418, 419, 420, 421, 422, 423, 424, 425, 426, 429, 430, 431, 432, 433, 434, 435, 436, 437, 438, 439, 440, 441, 442, 444, 445, 449, 450, 451, 453, 454, 455, 457, 458, 459, 460, 461, 462, 464, 466, 469, 470, 471, 472, 473, 474, 475, 476, 477, 478, 480, 481, 482, 483, 484, 485, 489, 490, 492, 493, 494, 495, 496, 497, 498, 499, 500, 501, 502, 504, 505, 506, 507, 508, 509, 510, 511, 514, 515, 516, 517, 518, 519, 521, 522, 523, 524, 525, 526, 527, 528, 530, 531, 533, 534, 535, 536, 537, 538, 539, 540, 541, 542, 543, 544, 545, 546, 549, 550, 551, 552, 553, 556, 557, 559, 560, 561, 562, 564, 565, 566, 567, 568, 569, 570, 571, 572, 573, 574, 575, 576, 577, 578, 579, 580, 581, 582, 583, 584, 585, 586, 587, 588, 589, 590, 592, 595, 596, 597, 598, 600, 601, 602, 603, 604, 605, 606, 607, 609, 610, 611, 614, 616, 617, 618, 619, 621, 622, 623, 624, 625, 627, 628, 631, 632, 633, 634, 636, 637, 638, 640, 643, 644, 645, 646, 647, 648, 649, 650, 652, 654, 655, 657, 660, 662, 663, 664, 665, 666, 667, 668, 669, 670, 671, 672, 673, 674, 675, 676, 677, 679, 680, 682, 684, 685, 686, 688, 689, 691, 692, 693, 694, 695, 696, 698, 699, 700, 701, 702, 704, 705, 706, 725, 729, 732, 735, 736, 737, 740, 741, 743, 744, 745, 746, 748, 749, 754, 755, 756, 763, 764, 766, 767, 768, 769, 770, 771, 772, 773, 774, 775, 776, 777, 778, 779, 780, 781, 782, 783, 786, 788, 789, 790, 791, 792, 793, 794, 795, 796, 798, 799, 800, 801, 802, 803, 804, 805, 806, 808, 810, 812, 813, 814, 815, 816, 817, 818, 819, 820, 821, 822, 823, 824, 825, 826, 827, 828, 829, 830, 831, 832, 833, 834, 836, 837, 838, 839, 840, 841, 842, 843, 844, 845, 846, 847, 848, 849, 850, 851, 852, 853, 854, 855, 856, 857, 858, 859, 860, 861, 862, 863, 864, 865, 866, 867, 868, 869, 870, 871, 872, 873, 874, 875, 876, 877, 878, 879, 880, 881, 882, 883, 884, 885, 886, 887, 888, 889, 890, 891, 892, 893, 894, 895, 896, 897, 898, 899, 901, 902, 903, 904, 905, 906, 907, 908, 910, 913, 914, 915, 916, 917, 918, 919, 920, 921, 922, 923, 924, 925, 926, 927, 928, 929, 930, 931, 932, 933, 934, 935, 936, 937, 938, 939, 940, 941, 942, 943, 944, 945, 946, 947, 948, 949, 950, 951, 953, 955, 956, 957, 958, 959, 960, 961, 962, 963, 964, 965, 966, 967, 968, 969, 970, 971, 972, 973, 974, 975, 977, 978, 979, 980, 981, 982, 983, 984, 985, 986, 987, 988, 989, 990, 991, 992, 993, 994, 995, 996, 997, 998, 999, 1000, 1001, 1002, 1003, 1004, 1005, 1006, 1007, 1008, 1009, 1010, 1011, 1012, 1013, 1014, 1015, 1016, 1017, 1018, 1019, 1020, 1021, 1022, 1023, 1024, 1025, 1026, 1027, 1028, 1029, 1030, 1031, 1032, 1033, 1034, 1035, 1036, 1037, 1038, 1039, 1041, 1042, 1043, 1044, 1046, 1047, 1048, 1049, 1050, 1051, 1052, 1053, 1054, 1055, 1056, 1057, 1058, 1059, 1060, 1061, 1062, 1063, 1064, 1065, 1066, 1067, 1069, 1070, 1071, 1072, 1073, 1074, 1075, 1076, 1077, 1078, 1079, 1080, 1081, 1082, 1083, 1084, 1085, 1086, 1087, 1088, 1089, 1090, 1091, 1092, 1093, 1094, 1095, 1096, 1097, 1098, 1099, 1100, 1101, 1102, 1103, 1104, 1106, 1114, 1115, 1116, 1117, 1118, 1119, 1120, 1121, 1122, 1123, 1124, 1125, 1126, 1127, 1128, 1129, 1130, 1131, 1132, 1133, 1134, 1135, 1136, 1137, 1138, 1139, 1141, 1142, 1145, 1146]'
$ws.Range("I3").Value = '{''Numk'': 14, ''KPar'': 14, ''Bucket_index'': 500}'
$ws.Range("J3").Value = 11.92300814692862
$ws.Range("O3").Value = '[233, 234, 235, 236, 237, 238, 239, 240, 241, 242, 243, 244, 245, 246, 247, 249, 250, 251, 253, 254, 255, 256, 257, 258, 259, 260, 261, 262, 264, 265, 266, 267, 268, 269, 270, 271, 272, 275, 276, 278, 279, 282, 283, 284, 287, 290, 291, 292, 301, 304, 305, 310, 311, 313, 315, 321, 323, 324, 325, 327, 328, 329, 330, 331, 332, 333, 334, 335, 336, 337, 338, 339, 341, 351, 355, 357, 358, 365, 366, 367, 368, 369, 370, 371, 377, 380, 392, 394, 398, 400, 404, 411, 412, 413, 414, 417, 419, 428, 429, 435, 436, 445, 447, 449, 452, 453, 465, 467, 475, 483, 484, 485, 486, 488, 495, 505, 506, 507, 512, 513, 514, 519, 520, 523, 535, 539, 542, 544, 546, 547, 550, 552, 553, 555, 560, 562, 565, 567, 570, 571, 572, 573, 574, 575, 576, 577, 578, 579, 580, 581, 582, 584, 586, 587, 588, 589, 590, 591, 592, 593, 594, 595, 596, 597, 598, 599, 600, 601, 602, 603, 604, 605, 606, 607, 608, 609, 610, 611, 612, 613, 614, 615, 616, 617, 618, 619, 620, 621, 622, 623, 624, 625, 626, 627, 628, 629, 630, 631, 632, 633, 634, 635, 636, 637, 638, 639, 640, 641, 642, 643, 644, 645, 646, 647, 648, 650, 651, 652, 653, 659, 660, 662, 663, 665, 667, 668, 674, 682, 683, 684, 686, 687, 690, 693, 698, 699, 701, 707, 708, 709, 710, 711, 712, 713, 714, 715, 716, 717, 718, 719, 720, 721, 722, 723, 724, 725, 726, 727, 728, 729, 730, 731, 732, 733, 734, 735, 736, 737, 738, 739, 740, 741, 742, 743, 744, 745, 746, 747, 748, 749, 750, 751, 752, 753, 754, 755, 756, 757, 758, 759, 760, 761, 762, 769, 771, 772, 775, 776, 778, 780, 781, 782, 785, 786, 787, 788, 790, 792, 793, 795, 797, 798, 799, 802, 803, 805, 806, 807, 808, 809, 811, 815, 816, 818, 819, 820, 821, 823, 824, 825, 826, 827, 828, 829, 830, 831, 832, 833, 834, 835, 836, 837, 838, 839, 840, 841, 842, 843, 844, 845, 846, 847, 848, 849, 850, 851, 852, 853, 854, 855, 856, 857, 858, 859, 860, 861, 862, 863, 864, 865, 866, 867, 868, 869, 870, 871, 872, 873, 874, 875, 876, 877, 878, 879, 880, 881, 882, 883, 884, 885, 886, 887, 888, 889, 890, 891, 892, 893, 894, 895, 896, 897, 898, 899, 900, 901, 902, 903, 904, 905, 906, 907, 908, 909, 910, 911, 912, 913, 914, 915, 916, 917, 918, 919, 920, 921, 922, 923, 924, 925, 926, 927, 928, 929, 930, 931, 932, 933, 934, 935, 936, 937, 938, 939, 940, 941, 942, 943, 944, 945, 946, 947, 948, 949, 950, 951, 952, 953, 954, 955, 956, 957, 958, 959, 960, 961, 962, 963, 964, 965, 966, 967, 968, 969, 970, 971, 972, 973, 974, 975, 976, 977, 978, 979, 980, 981, 982, 983, 984, 985, 987, 988, 989, 990, 991, 992, 993, 994, 995, 996, 997, 998, 999, 1000, 1001, 1002, 1003, 1004, 1005, 1006, 1007, 1008, 1009, 1010, 1011, 1012, 1013, 1014, 1015, 1016, 1017, 1018, 1019, 1020, 1021, 1022, 1023, 1024, 1025, 1026, 1027, 1028, 1029, 1030, 1031, 1032, 1033, 1034, 1035, 1036, 1037, 1038, 1039, 1040, 1041, 1042, 1043, 1044, 1045, 1046, 1047, 1048, 1049, 1050, 1051, 1052, 1053, 1054, 1055, 1056, 1057, 1058, 1059, 1060, 1061, 1062, 1063, 1064, 1065, 1066, 1067, 1068, 1069, 1070, 1071, 1072, 1073, 1074, 1075, 1076, 1077, 1078, 1079, 1080, 1081, 1082, 1083, 1084, 1085, 1086, 1087, 1088, 1089, 1090, 1091, 1092, 1093, 1094, 1095, 1096, 1097, 1098, 1099, 1100, 1101, 1102, 1103, 1104, 1105, 1106, 1107, 1108, 1109, 1110, 1111, 1112, 1113, 1114, 1115, 1116, 1117, 1118, 1119, 1120, 1121, 1122, 1123, 1124, 1125, 1126, 1127, 1128, 1129, 1130, 1131, 1132, 1133, 1134, 1135, 1136, 1137, 1138, 1139, 1140, 1141, 1142, 1143, 1144, 1145, 1146]'
$ws.Range("P3").Value = 0.1213872832369942
$ws.Range("Q3").Value = '{''window_size'': 232.0, ''max_size_ae'': 2}'
$ws.Range("R3").Value = 0.4022495881654322
$ws.Range("S3").Value = '[1, 2, 3, 4, 5, 6, 9, 10, 11, 14, 20, 22, 26, 32, 34, 37, 38, 39, 41, 43, 45, 47, 50, 51, 52, 54, 55, 56, 60, 69, 70, 84, 101, 102, 107, 109, 112, 114, 119, 121, 122, 128, 133, 134, 138, 139, 144, 147, 148, 151, 155, 166, 174, 175, 176, 180, 193, 196, 198, 200, 204, 212, 216, 217, 220, 223, 224, 225, 230, 238, 239, 240, 241, 244, 247, 249, 251, 254, 258, 267, 268, 269, 271, 273, 278, 287, 298, 299, 300, 301, 311, 313, 325, 327, 330, 331, 332, 337, 340, 344, 351, 352, 355, 357, 358, 363, 365, 366, 367, 368, 370, 372, 374, 377, 380, 381, 382, 386, 390, 391, 394, 395, 396, 399, 400, 401, 402, 403, 404, 406, 413, 414, 417, 419, 429, 436, 441, 445, 449, 456, 457, 459, 461, 462, 463, 464, 465, 470, 472, 475, 476, 483, 485, 490, 495, 496, 497, 499, 500, 505, 506, 513, 514, 518, 519, 523, 526, 531, 536, 538, 539, 540, 542, 543, 549, 550, 552, 553, 554, 555, 556, 560, 561, 562, 564, 565, 566, 570, 572, 573, 574, 575, 576, 579, 580, 581, 582, 584, 585, 587, 590, 591, 592, 595, 597, 598, 599, 600, 601, 602, 603, 604, 605, 606, 607, 608, 609, 610, 611, 612, 613, 614, 615, 616, 617, 618, 619, 620, 621, 622, 623, 624, 625, 626, 627, 628, 629, 630, 631, 632, 633, 634, 636, 637, 638, 639, 640, 641, 643, 644, 645, 647, 648, 650, 652, 653, 657, 658, 659, 660, 662, 663, 665, 667, 668, 671, 674, 676, 677, 679, 681, 682, 684, 689, 692, 693, 698, 699, 700, 701, 702, 703, 707, 709, 710, 711, 712, 713, 714, 715, 716, 717, 718, 719, 720, 721, 722, 723, 724, 725, 726, 727, 728, 729, 730, 731, 732, 733, 734, 736, 737, 738, 740, 741, 742, 743, 744, 745, 746, 747, 749, 750, 751, 752, 753, 754, 755, 756, 758, 759, 760, 761, 762, 763, 765, 766, 769, 772, 773, 775, 776, 777, 778, 780, 781, 782, 783, 784, 786, 787, 788, 789, 792, 793, 794, 795, 796, 797, 798, 799, 800, 802, 803, 804, 805, 806, 807, 808, 809, 810, 811, 812, 813, 815, 816, 817, 818, 819, 820, 821, 822, 823, 824, 825, 826, 827, 828, 829, 830, 831, 832, 833, 834, 835, 836, 837, 838, 839, 840, 841, 842, 843, 844, 845, 846, 847, 848, 849, 850, 851, 852, 853, 854, 855, 856, 857, 858, 859, 860, 861, 862, 863, 864, 865, 866, 867, 868, 869, 870, 871, 872, 873, 874, 875, 876, 877, 878, 879, 880, 881, 882, 883, 884, 885, 886, 887, 888, 889, 890, 891, 892, 893, 894, 895, 896, 897, 898, 899, 900, 901, 902, 903, 904, 905, 906, 907, 908, 909, 910, 911, 912, 913, 914, 915, 916, 917, 918, 919, 920, 921, 922, 923, 924, 925, 926, 927, 928, 929, 930, 931, 932, 933, 934, 935, 936, 937, 938, 939, 940, 941, 942, 943, 944, 945, 946, 947, 948, 949, 950, 951, 952, 953, 954, 955, 956, 957, 958, 959, 960, 961, 962, 963, 965, 966, 967, 968, 969, 970, 971, 972, 973, 975, 976, 977, 978, 979, 980, 981, 983, 984, 985, 987, 988, 990, 991, 993, 1000, 1004, 1005, 1013, 1030, 1035, 1036, 1037, 1040, 1044, 1061, 1063, 1068, 1077, 1082, 1083, 1085, 1089, 1091, 1097, 1100, 1101, 1106, 1107, 1109, 1110, 1111, 1112, 1113, 1115, 1118, 1120, 1130, 1141, 1142, 1143, 1144, 1146]'
$ws.Range("T3").Value = 0.1221864951768489
$ws.Range("U3").Value = '{''max_features'': 7, ''window_size'': 245, ''n_estimator'': 40}'
$ws.Range("V3").Value = 31.77612683898769

$ws.Range("G4").Value = '[3, 37, 200, 217, 227, 384, 396, 408, 578, 657, 799]'
$ws.Range("I4").Value = '{''Numk'': 13, ''KPar'': 4, ''Bucket_index'': 500}'
$ws.Range("J4").Value = 5.977500962093472
$ws.Range("O4").Value = '[179, 906, 907, 908, 911, 914, 915, 916, 917, 918, 919, 920, 921, 922, 923, 924, 925, 926, 927, 928, 929, 930, 931, 932, 933, 934, 935, 936, 937, 940, 942, 994, 995, 996, 997, 998, 999, 1000, 1001, 1002, 1003, 1004, 1005, 1006, 1007, 1008, 1009, 1010, 1011, 1012, 1013, 1014, 1015, 1016, 1017, 1018, 1019, 1020, 1021, 1022, 1023, 1024, 1025, 1026, 1027, 1028, 1029, 1030, 1031, 1032, 1033, 1034, 1035, 1036, 1037, 1038, 1039, 1040, 1041, 1042, 1043, 1044, 1045, 1046, 1047, 1048, 1049, 1050, 1051, 1052, 1053, 1054, 1055, 1056, 1057, 1058, 1059, 1060, 1061, 1062, 1063, 1064, 1065, 1066, 1067, 1068, 1069, 1070, 1071, 1072, 1073, 1074, 1075, 1076, 1077, 1078, 1079, 1080, 1081, 1082, 1083, 1084, 1085, 1086, 1087, 1088, 1089, 1090, 1091, 1092, 1093, 1094, 1095, 1096, 1097, 1098, 1099, 1100, 1101, 1102, 1103, 1104, 1105, 1106, 1107, 1108, 1109, 1110, 1111, 1112, 1113, 1114, 1115, 1116, 1117, 1118, 1119, 1120, 1121, 1122, 1123, 1124, 1125, 1126, 1127, 1128, 1129, 1130, 1131, 1132, 1133, 1134, 1135, 1136, 1137, 1138, 1139, 1140, 1141, 1142, 1143, 1144, 1145, 1146, 1147, 1148, 1149, 1150, 1151, 1152, 1153, 1154, 1155, 1156, 1157, 1158, 1159, 1160, 1161, 1162, 1163, 1164, 1165, 1166, 1167, 1168, 1169, 1170, 1171, 1172, 1173, 1174, 1175, 1176, 1177, 1178, 1179, 1180, 1181, 1182, 1183, 1184, 1185, 1186, 1187, 1188, 1189, 1190]'
$ws.Range("P4").Value = 0.2784810126582278
$ws.Range("Q4").Value = '{''window_size'': 178.0, ''max_size_ae'': 3}'
$ws.Range("R4").Value = 0.2346447729505599
$ws.Range("S4").Value = '[208, 233, 244, 245, 246, 248, 249, 250, 266, 285, 347, 373, 422, 424, 433, 437, 444, 449, 454, 455, 469, 472, 494, 516, 519, 580, 582, 597, 607, 608, 616, 633, 635, 647, 648, 650, 694, 695, 705, 716, 719, 737, 767, 768, 771, 774, 782, 784, 785, 786, 788, 789, 790, 791, 792, 793, 794, 795, 796, 797, 798, 800, 801, 803, 805, 806, 809, 811, 814, 817, 822, 823, 829, 832, 833, 836, 842, 843, 844, 845, 847, 848, 850, 851, 852, 853, 854, 855, 856, 857, 858, 859, 860, 861, 862, 863, 864, 865, 866, 867, 868, 869, 870, 871, 872, 873, 874, 875, 876, 877, 878, 879, 880, 881, 882, 883, 884, 885, 886, 887, 888, 889, 890, 891, 892, 894, 895, 897, 899, 905, 906, 907, 908, 909, 912, 914, 916, 917, 921, 922, 924, 926, 930, 931, 932, 933, 934, 935, 936, 938, 939, 940, 941, 942, 943, 944, 945, 946, 947, 948, 949, 950, 951, 952, 953, 954, 955, 956, 957, 958, 959, 960, 961, 962, 963, 964, 965, 966, 968, 969, 971, 972, 973, 974, 975, 976, 978, 979, 981, 985, 986, 988, 989, 990, 991, 992, 993, 994, 995, 996, 997, 998, 1001, 1002, 1003, 1004, 1006, 1007, 1008, 1011, 1013, 1015, 1016, 1017, 1023, 1025, 1026, 1027, 1028, 1029, 1030, 1031, 1033, 1034]'
$ws.Range("T4").Value = 0.2222222222222222
$ws.Range("U4").Value = '{''max_features'': 7, ''window_size'': 207, ''n_estimator'': 33}'
$ws.Range("V4").Value = 31.3619844308123

$ws.Range("G5").Value = '[563]'
$ws.Range("I5").Value = '{''Numk'': 13, ''KPar'': 9, ''Bucket_index'': 500}'
$ws.Range("J5").Value = 9.123951609013602
$ws.Range("O5").Value = '[147, 547, 549, 550, 551, 552, 554, 555, 556, 557, 558, 559, 560, 561, 562, 563, 564, 565, 566, 567, 568, 569, 570, 571, 572, 573, 574, 575, 576, 577, 578, 580, 581, 582, 583, 584, 585, 586, 587, 588, 589, 590, 591, 592, 593, 594, 595, 596, 597, 598, 599, 600, 601, 602, 603, 604, 605, 606, 607, 608, 609, 610, 611, 612, 613, 614, 615, 616, 617, 618, 619, 620, 621, 622, 623, 624, 625, 626, 627, 628, 629, 630, 631, 632, 633, 634, 635, 636, 637, 638, 639, 640, 641, 642, 643, 644, 645, 646, 647, 648, 649, 650, 651, 652, 653, 654, 655, 656, 657, 658, 659, 660, 661, 662, 663, 664, 665, 666, 667, 668, 669, 670, 671, 672, 673, 674, 675, 676, 677, 678, 679, 680, 681, 682, 683, 684, 685, 686, 687, 688, 689, 690, 691, 692, 693, 694, 695, 696, 697, 698, 699, 700, 701, 702, 703, 704, 705, 706, 707, 708, 709, 710, 711, 712, 713, 714, 715, 716, 717, 718, 719, 720, 721, 722, 723, 724, 725, 726, 727, 728, 729, 730, 731, 732, 733, 745, 746, 747, 748, 749, 750]'
$ws.Range("P5").Value = 0.2477064220183486
$ws.Range("Q5").Value = '{''window_size'': 146.0, ''max_size_ae'': 6}'
$ws.Range("R5").Value = 0.2529439418576658
$ws.Range("S5").Value = '[2, 11, 17, 36, 50, 88, 94, 113, 169, 178, 181, 182, 184, 188, 189, 190, 191, 193, 194, 195, 196, 202, 204, 209, 212, 214, 216, 217, 218, 219, 220, 222, 223, 225, 230, 231, 232, 234, 235, 236, 239, 245, 246, 247, 248, 251, 252, 255, 257, 259, 261, 264, 267, 271, 274, 275, 277, 278, 282, 287, 291, 293, 296, 297, 298, 299, 300, 301, 302, 303, 304, 306, 307, 308, 309, 310, 311, 312, 313, 314, 315, 317, 318, 319, 320, 324, 325, 326, 327, 328, 330, 331, 333, 334, 335, 336, 337, 339, 341, 342, 344, 345, 346, 348, 350, 351, 352, 354, 355, 356, 357, 359, 360, 361, 362, 363, 364, 365, 367, 369, 370, 371, 372, 375, 378, 379, 381, 382, 386, 387, 388, 390, 391, 392, 393, 394, 395, 397, 398, 400, 401, 402, 403, 405, 406, 407, 408, 409, 411, 412, 413, 414, 415, 416, 419, 420, 421, 422, 423, 424, 425, 426, 427, 428, 429, 430, 431, 432, 433, 434, 435, 436, 437, 438, 439, 441, 442, 444, 445, 446, 447, 448, 449, 450, 451, 452, 454, 455, 457, 458, 459, 461, 462, 464, 465, 466, 467, 468, 469, 470, 471, 472, 474, 475, 477, 478, 479, 485, 486, 494, 500, 501, 502, 503, 506, 507, 508, 510, 511, 512, 513, 514, 515, 516, 517, 518, 519, 520, 521, 522, 523, 524, 526, 527, 529, 531, 532, 533, 534, 535, 536, 537, 538, 539, 540, 541, 542, 543, 544, 546, 547, 548, 549, 550, 551, 552, 553, 554, 556, 558, 559, 560, 561, 562, 563, 564, 566, 567, 568, 569, 570, 571, 572, 573, 574, 575, 576, 577, 578, 579, 580, 582, 584, 585, 586, 587, 588, 589, 591, 592, 593, 594, 595, 596, 597, 598, 599, 612, 620, 623, 627, 632, 645, 660, 661, 694, 726, 745, 746, 747, 748, 749, 750]'
$ws.Range("T5").Value = 0.08562691131498472
$ws.Range("U5").Value = '{''max_features'': 7, ''window_size'': 200, ''n_estimator'': 37}'
$ws.Range("V5").Value = 27.25725762988441

$ws.Range("G6").Value = '[3, 7, 15, 17, 20, 23, 24, 27, 28, 38, 45, 61, 64, 65, 67, 71, 75, 78, 83, 89, 96, 103, 105, 110, 116, 130, 133, 134, 139, 140, 143, 144, 146, 151, 156, 161, 162, 164, 165, 166, 176, 179, 180, 181, 182, 186, 201, 202, 205, 211, 228, 232, 234, 236, 240, 249, 250, 255, 260, 261, 267, 278, 279, 280, 281, 284, 286, 290, 291, 293, 300, 307, 310, 311, 316, 317, 321, 328, 331, 333, 334, 335, 340, 346, 347, 348, 351, 352, 362, 363, 364, 365, 366, 368, 369, 371, 376, 377, 379, 382, 385, 386, 387, 389, 390, 391, 393, 394, 395, 397, 398, 401, 403, 404, 406, 411, 412, 413, 417, 419, 421, 431, 434, 436, 438, 443, 445, 447, 449, 451, 452, 455, 461, 466, 470, 471, 474, 475, 480, 481, 495, 496, 497, 498, 505, 508, 510, 513, 516, 525, 528, 533, 534, 535, 536, 541, 543, 545, 549, 551, 556, 559, 560, 561, 563, 564, 565, 566, 567, 568, 569, 570, 571, 572, 573, 574, 575, 576, 578, 579, 580, 581, 582, 586, 588, 591, 597, 598, 599, 603, 609, 610, 611, 612, 615, 617, 619, 620, 624, 628, 631, 632, 633, 634, 635, 638, 641, 642, 643, 644, 645, 648, 650, 652, 655, 656, 658, 660, 663, 664, 665, 666, 667, 671, 674, 675, 676, 678, 681, 686, 689, 690, 693, 705, 706, 710, 713, 717, 724, 732, 738, 739, 743, 747, 764, 765, 767, 769, 770, 774, 777, 778, 783, 804, 807, 813, 816, 819, 827, 828, 830, 831, 834, 835, 836, 837, 842, 843, 846, 850, 851, 855, 858, 860, 862, 865, 867, 882, 884, 905, 906, 912, 920, 921, 933, 936, 942, 949, 956, 957, 958, 959, 970, 972, 973, 976, 980, 981, 986, 989, 992, 993, 1000, 1002, 1004, 1005, 1007, 1011, 1012, 1014, 1017, 1018, 1019, 1022, 1023, 1028, 1029, 1030, 1031, 1035, 1037, 1038, 1039, 1041, 1094, 1099, 1100, 1106, 1112, 1113, 1114, 1119, 1121, 1124, 1127, 1131, 1132, 1136, 1137, 1139, 1140, 1143, 1146, 1148, 1149, 1150]'
$ws.Range("I6").Value = '{''Numk'': 12, ''KPar'': 13, ''Bucket_index'': 500}'
$ws.Range("J6").Value = 21.83030613115989
$ws.Range("O6").Value = '[267, 579]'
$ws.Range("P6").Value = 0.5
$ws.Range("Q6").Value = '{''window_size'': 266.0, ''max_size_ae'': 5}'
$ws.Range("R6").Value = 0.4745001432020217
$ws.Range("S6").Value = '[991]'
$ws.Range("T6").Value = 0.6666666666666666
$ws.Range("U6").Value = '{''max_features'': 7, ''window_size'': 210, ''n_estimator'': 35}'
$ws.Range("V6").Value = 36.63673501391895

$ws.Range("G7").Value = '[5, 6, 9, 11, 16, 17, 28, 35, 39, 42, 44, 46, 48, 51, 54, 56, 58, 63, 64, 69, 70, 74, 84, 86, 92, 93, 104, 105, 106, 107, 108, 117, 118, 119, 121, 124, 126, 129, 131, 132, 133, 134, 135, 137, 139, 141, 143, 145, 146, 147, 148, 149, 150, 151, 152, 153, 154, 158, 160, 163, 165, 168, 170, 172, 183, 189, 197, 201, 202, 212, 218, 229, 230, 240, 241, 243, 245, 255, 261, 264, 272, 278, 281, 283, 284, 285, 286, 287, 288, 293, 302, 307, 311, 315, 316, 317, 319, 320, 321, 325, 326, 331, 337, 341, 344, 350, 357, 363, 371, 375, 378, 381, 382, 384, 386, 387, 389, 401, 402, 403, 404, 405, 409, 410, 411, 416, 417, 420, 423, 427, 430, 431, 432, 433, 434, 437, 438, 440, 444, 449, 452, 455, 457, 460, 461, 465, 466, 467, 471, 472, 476, 479, 490, 495, 499, 503, 505, 513, 518, 519, 520, 527, 528, 529, 530, 531, 533, 536, 537, 538, 540, 541, 543, 545, 550, 552, 553, 557, 563, 565, 566, 569, 570, 571, 572, 573, 574, 575, 576, 579, 582, 583, 585, 589, 592, 593, 594, 596, 601, 612, 615, 617, 618, 623, 627, 630, 633, 634, 635, 636, 638, 639, 640, 643, 644, 648, 651, 652, 654, 655, 656, 658, 659, 663, 664]'
$ws.Range("I7").Value = '{''Numk'': 7, ''KPar'': 13, ''Bucket_index'': 500}'
$ws.Range("J7").Value = 6.494259615894407
$ws.Range("O7").Value = '[145, 146, 395, 492]'
$ws.Range("P7").Value = 0.6
$ws.Range("Q7").Value = '{''window_size'': 144.0, ''max_size_ae'': 7}'
$ws.Range("R7").Value = 0.1674506790004671
$ws.Range("S7").Value = '[104, 105, 106, 107, 108, 128, 250, 297, 298, 483, 488, 523]'
$ws.Range("T7").Value = 0.5
$ws.Range("U7").Value = '{''max_features'': 7, ''window_size'': 200, ''n_estimator'': 19}'
$ws.Range("V7").Value = 30.90919523593038

$ws.Range("G8").Value = '[19, 25, 30, 56, 72, 75, 129, 180, 184, 220, 231, 262, 284, 338, 375, 383, 450, 454, 503, 538, 548, 568, 582, 583, 606, 641, 667, 734, 759, 779, 783, 797, 861, 874, 885, 899, 905, 913, 914, 930, 963, 992, 1007, 1010, 1029, 1034]'
$ws.Range("I8").Value = '{''Numk'': 8, ''KPar'': 3, ''Bucket_index'': 500}'
$ws.Range("J8").Value = 4.80887848418206
$ws.Range("O8").Value = '[188, 497, 524, 526, 533, 534, 536, 537, 538, 539, 542, 543, 544, 546, 547, 548, 549, 550, 551, 552, 553, 554, 555, 556, 557, 558, 559, 560, 561, 562, 563, 564, 565, 566, 567, 568, 569, 570, 571, 572, 573, 574, 575, 576, 577, 578, 579, 580, 581, 582, 583, 584, 585, 586, 587, 588, 589, 590, 591, 592, 593, 594, 595, 596, 597, 598, 599, 600, 601, 602, 603, 604, 605, 606, 607, 608, 609, 610, 611, 612, 613, 614, 615, 616, 617, 618, 619, 620, 621, 622, 623, 624, 625, 626, 627, 628, 629, 630, 631, 632, 633, 634, 635, 636, 637, 638, 639, 640, 641, 642, 643, 644, 645, 646, 647, 648, 649, 650, 651, 652, 653, 654, 655, 656, 657, 658, 659, 660, 661, 662, 663, 664, 665, 666, 667, 668, 669, 670, 671, 672, 673, 674, 675, 676, 677, 678, 679, 680, 681, 682, 683, 684, 685, 686, 687, 688, 689, 690, 691, 692, 693, 694, 695, 696, 697, 698, 699, 700, 701, 702, 703, 704, 705, 706, 707, 708, 709, 710, 711, 712, 713, 714, 715, 716, 717, 718, 719, 720, 721, 722, 723, 724, 725, 726, 727, 728, 729, 730, 731, 732, 733, 734, 735, 736, 737, 738, 739, 740, 741, 742, 743, 744, 745, 746, 747, 748, 749, 750, 751, 752, 753, 754, 755, 756, 757, 758, 759, 760, 761, 762, 763, 764, 765, 766, 767, 768, 769, 770, 771, 772, 773, 774, 775, 776, 777, 778, 779, 780, 781, 782, 836, 837, 838, 839, 840, 841, 842, 843, 844, 845, 846, 847, 848, 849, 850, 851, 852, 853, 854, 855, 856, 857, 858, 859, 860, 861, 862, 863, 864, 865, 866, 867, 868, 869, 870, 871, 872, 873, 874, 875, 876, 877, 878, 879, 880, 881, 882, 883, 884, 885, 886, 887, 888, 889, 890, 891, 892, 893, 894, 895, 896, 897, 898, 899, 900, 901, 902, 903, 904, 905, 906, 907, 908, 909, 910, 911, 912, 913, 914, 915, 916, 917, 918, 919, 920, 921, 922, 923, 924, 925, 926, 927, 928, 929, 930, 931, 932, 933, 934, 935, 936, 937, 938, 939, 940, 941, 942, 943, 944, 945, 946, 947, 948, 949, 950, 951, 952, 953, 954, 955, 956, 957, 958, 959, 960, 961, 962, 963, 964, 965, 966, 967, 968, 969, 970, 971, 972, 973, 974, 975, 976, 977, 978, 979, 980, 981, 982, 983, 984, 985, 986, 987, 988, 989, 990, 991, 992, 993, 994, 995, 996, 997, 998, 999, 1000, 1001, 1002, 1003, 1004, 1005, 1006, 1007, 1008, 1009, 1010, 1011, 1012, 1013, 1014, 1015, 1016, 1017, 1018, 1019, 1020, 1021, 1022, 1023, 1024, 1025, 1026, 1027, 1028, 1029, 1030, 1031, 1032, 1033, 1034, 1035, 1036, 1037, 1038, 1039, 1040, 1041, 1042, 1043, 1044, 1045, 1046, 1047, 1048, 1049, 1050, 1051, 1052, 1053, 1054, 1055, 1056, 1057, 1058, 1059, 1060, 1061, 1062, 1063, 1064, 1065]'
$ws.Range("P8").Value = 0.2857142857142857
$ws.Range("Q8").Value = '{''window_size'': 187.0, ''max_size_ae'': 4}'
$ws.Range("R8").Value = 0.174202990019694
$ws.Range("S8").Value = '[1, 237, 248, 268, 324, 397, 403, 409, 410, 415, 440, 497, 498, 501, 502, 515, 544, 558, 561, 580, 585, 651, 652, 662, 688, 721, 727, 730, 732, 764, 767, 773, 799, 809, 827, 828, 834, 840, 848, 849, 850, 851, 852, 854, 857, 876, 902, 1066, 1070]'
$ws.Range("T8").Value = 0.2950819672131148
$ws.Range("U8").Value = '{''max_features'': 7, ''window_size'': 228, ''n_estimator'': 20}'
$ws.Range("V8").Value = 30.60310626681894

$ws.Range("G9").Value = '[2, 9, 18, 19, 39, 42, 45, 48, 50, 52, 53, 55, 61, 62, 74, 81, 89, 99, 119, 120, 123, 127, 139, 173, 197, 201, 202, 213, 228, 232, 234, 253, 262, 266, 267, 272, 283, 289, 307, 313, 315, 319, 327, 333, 335, 339, 340, 343, 344, 355, 391, 394, 398, 400, 410, 419, 421, 423, 424, 432, 433, 438, 441, 443, 451, 461, 464, 465, 473, 474, 476, 477, 478, 481, 496, 501, 503, 505, 513, 526, 528, 535, 536, 559, 561, 562, 564, 567, 579, 583, 584, 586, 587, 591, 605, 607, 614, 621, 622, 629, 631, 671, 676, 677, 685, 692, 697, 734, 763, 785, 802, 805, 806, 810, 811, 812, 821, 827, 857, 906, 913, 915, 920, 921, 949, 956, 1002, 1005, 1015, 1021, 1053, 1059, 1061, 1062, 1092, 1095, 1119, 1130, 1134]'
$ws.Range("I9").Value = '{''Numk'': 9, ''KPar'': 7, ''Bucket_index'': 500}'
$ws.Range("J9").Value = 7.070770997088403
$ws.Range("O9").Value = '[200, 847, 848, 849, 850, 851, 852, 853, 858, 859, 860, 861, 862, 863, 864, 865, 866, 868, 874, 875, 876, 878, 885, 887, 888, 889, 890, 891, 892, 894, 896, 898, 905, 906, 908, 910, 916, 918, 920, 922, 923, 924, 934, 936, 948, 950, 952, 953, 965, 967, 969, 971, 977, 979, 981, 983, 997, 999, 1012, 1014, 1016, 1027, 1038, 1042, 1052, 1054, 1055, 1064, 1066, 1069, 1077, 1078, 1079, 1080, 1082, 1083, 1084, 1086, 1093, 1095, 1097, 1099, 1108, 1109, 1111, 1112, 1113]'
$ws.Range("P9").Value = 0.2926829268292683
$ws.Range("Q9").Value = '{''window_size'': 199.0, ''max_size_ae'': 1}'
$ws.Range("R9").Value = 0.385773902060464
$ws.Range("S9").Value = '[888, 892, 919, 984, 1146]'
$ws.Range("T9").Value = 0.4444444444444445
$ws.Range("U9").Value = '{''max_features'': 7, ''window_size'': 202, ''n_estimator'': 25}'
$ws.Range("V9").Value = 32.27393499505706
